function Set-CellText($ws, $ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText $ws "D2" "68.530.45"
Set-CellText $ws "E2" "  -0.45%  "
Set-CellText $ws "D3" "2.428.77"
Set-CellText $ws "E3" "  -1.36%  "
Set-CellText $ws "E4" "  -0.07%  "
Set-CellText $ws "D5" "558.61"
Set-CellText $ws "E5" "  -0.04%  "
Set-CellText $ws "D6" "160.55"
Set-CellText $ws "E6" "  -0.64%  "
Set-CellText $ws "E7" "  -0.05%  "
Set-CellText $ws "D8" "0.509"
Set-CellText $ws "E8" "  +0.72%  "
Set-CellText $ws "E9" "  +9.20%  "
Set-CellText $ws "E10" "  -1.60%  "
Set-CellText $ws "E11" "  +0.14%  "
Set-CellText $ws "E12" "  -5.12%  "
Set-CellText $ws "D13" "68.407.84"
Set-CellText $ws "E13" "  -0.52%  "
Set-CellText $ws "D14" "2.874.03"
Set-CellText $ws "E14" "  -1.09%  "
Set-CellText $ws "E15" "  +3.56%  "
Set-CellText $ws "D16" "23.09"
Set-CellText $ws "E16" "  -1.88%  "
Set-CellText $ws "D17" "2.426.37"
Set-CellText $ws "E17" "  -1.14%  "
Set-CellText $ws "D18" "10.48"
Set-CellText $ws "E18" "  -1.89%  "
Set-CellText $ws "D19" "335.37"
Set-CellText $ws "E19" "  +0.09%  "
Set-CellText $ws "E20" "  -0.82%  "
Set-CellText $ws "E21" "  +1.14%  "
Set-CellText $ws "E22" "  +2.40%  "
Set-CellText $ws "E23" "  -0.03%  "
Set-CellText $ws "D24" "66.77"
Set-CellText $ws "E24" "  -0.12%  "
Set-CellText $ws "D25" "3.67"
Set-CellText $ws "E25" "  +0.66%  "
Set-CellText $ws "D26" "2.553.30"
Set-CellText $ws "E26" "  -1.83%  "
Set-CellText $ws "D27" "0.999"
Set-CellText $ws "E27" "  -0.19%  "
Set-CellText $ws "D28" "8.21"
Set-CellText $ws "E28" "  +0.79%  "
Set-CellText $ws "D29" "0.0₃0820"
Set-CellText $ws "E29" "  +0.94%  "
Set-CellText $ws "D30" "7.13"
Set-CellText $ws "E30" "  -0.41%  "
Set-CellText $ws "E31" "  -0.08%  "
Set-CellText $ws "D32" "427.37"
Set-CellText $ws "E32" "  -0.02%  "
Set-CellText $ws "E33" "  +0.75%  "
Set-CellText $ws "E34" "  -0.29%  "
Set-CellText $ws "D35" "159.77"
Set-CellText $ws "E35" "  +1.18%  "
Set-CellText $ws "D36" "19.05"
Set-CellText $ws "E36" "  +0.13%  "
Set-CellText $ws "E38" "  +0.87%  "
Set-CellText $ws "E39" "  -3.10%  "
Set-CellText $ws "E40" "  -0.42%  "
Set-CellText $ws "E41" "  -1.44%  "
Set-CellText $ws "E42" "  +2.16%  "
Set-CellText $ws "E43" "  +1.02%  "
Set-CellText $ws "D44" "2.05"
Set-CellText $ws "E44" "  +0.32%  "
Set-CellText $ws "D45" "131.36"
Set-CellText $ws "E45" "  -0.36%  "
Set-CellText $ws "E46" "  +0.08%  "
Set-CellText $ws "E47" "  -0.08%  "
Set-CellText $ws "E48" "  -0.06%  "
Set-CellText $ws "E49" "  -0.63%  "
Set-CellText $ws "D50" "0.0915"
Set-CellText $ws "E50" "  +0.49%  "
Set-CellText $ws "E51" "  +0.26%  "
